$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the updated cells remain stored as text (matching the original inline-string
# representation of numbers/percentages) instead of being auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"
$ws.Range("G2:G51").NumberFormat = "@"

$ws.Range("D2").Value = "301.86"; $ws.Range("E2").Value = "-0.92%"; $ws.Range("G2").Value = "7"
$ws.Range("D3").Value = "37.54"; $ws.Range("E3").Value = "7.52%"; $ws.Range("G3").Value = "7"
$ws.Range("D4").Value = "5.010"; $ws.Range("E4").Value = "-2.64%"; $ws.Range("G4").Value = "7"
$ws.Range("D5").Value = "0.07857"; $ws.Range("E5").Value = "1.18%"; $ws.Range("G5").Value = "7"
$ws.Range("D6").Value = "2.235"; $ws.Range("E6").Value = "-5.07%"; $ws.Range("G6").Value = "7"
$ws.Range("D7").Value = "8.033"; $ws.Range("E7").Value = "0.18%"; $ws.Range("G7").Value = "7"
$ws.Range("D8").Value = "4.020"; $ws.Range("E8").Value = "2.09%"; $ws.Range("G8").Value = "7"
$ws.Range("D9").Value = "0.9088"; $ws.Range("E9").Value = "-1.43%"; $ws.Range("G9").Value = "7"
$ws.Range("D10").Value = "0.1882"; $ws.Range("E10").Value = "4.60%"; $ws.Range("G10").Value = "7"
$ws.Range("D11").Value = "0.09266"; $ws.Range("E11").Value = "-7.38%"; $ws.Range("G11").Value = "7"
$ws.Range("D12").Value = "0.08488"; $ws.Range("E12").Value = "-0.69%"; $ws.Range("G12").Value = "7"
$ws.Range("D13").Value = "0.03524"; $ws.Range("E13").Value = "6.65%"; $ws.Range("G13").Value = "7"
$ws.Range("D14").Value = "0.09967"; $ws.Range("E14").Value = "0.77%"; $ws.Range("G14").Value = "7"
$ws.Range("D15").Value = "0.001484"; $ws.Range("E15").Value = "-0.22%"; $ws.Range("G15").Value = "7"
$ws.Range("D16").Value = "0.005701"; $ws.Range("E16").Value = "-1.42%"; $ws.Range("G16").Value = "7"
$ws.Range("D17").Value = "3.463"; $ws.Range("E17").Value = "-0.18%"; $ws.Range("G17").Value = "7"
$ws.Range("D18").Value = "2.063"; $ws.Range("E18").Value = "-3.25%"; $ws.Range("G18").Value = "7"
$ws.Range("G19").Value = "7"
$ws.Range("D20").Value = "0.1308"; $ws.Range("E20").Value = "0.51%"; $ws.Range("G20").Value = "7"
$ws.Range("D21").Value = "4.759"; $ws.Range("E21").Value = "10.60%"; $ws.Range("G21").Value = "7"
$ws.Range("D22").Value = "0.2204"; $ws.Range("E22").Value = "-7.52%"; $ws.Range("G22").Value = "7"
$ws.Range("D23").Value = "0.04652"; $ws.Range("E23").Value = "1.95%"; $ws.Range("G23").Value = "7"
$ws.Range("D24").Value = "0.001230"; $ws.Range("E24").Value = "1.15%"; $ws.Range("G24").Value = "7"
$ws.Range("D25").Value = "0.004452"; $ws.Range("E25").Value = "-0.19%"; $ws.Range("G25").Value = "7"
$ws.Range("D26").Value = "0.0001300"; $ws.Range("E26").Value = "0.18%"; $ws.Range("G26").Value = "7"
$ws.Range("D27").Value = "0.0004751"; $ws.Range("E27").Value = "28.60%"; $ws.Range("G27").Value = "7"
$ws.Range("G28").Value = "7"
$ws.Range("G29").Value = "7"
$ws.Range("G30").Value = "7"
$ws.Range("G31").Value = "7"
$ws.Range("G32").Value = "7"
$ws.Range("G33").Value = "7"
$ws.Range("G34").Value = "7"
$ws.Range("G35").Value = "7"
$ws.Range("G36").Value = "7"
$ws.Range("G37").Value = "7"
$ws.Range("G38").Value = "7"
$ws.Range("E39").Value = "-1.28%"; $ws.Range("G39").Value = "7"
$ws.Range("D40").Value = "0.04746"; $ws.Range("E40").Value = "-0.25%"; $ws.Range("G40").Value = "7"
$ws.Range("D41").Value = "0.007824"; $ws.Range("E41").Value = "0.74%"; $ws.Range("G41").Value = "7"
$ws.Range("E42").Value = "-1.40%"; $ws.Range("G42").Value = "7"
$ws.Range("D43").Value = "0.007662"; $ws.Range("E43").Value = "8.16%"; $ws.Range("G43").Value = "7"
$ws.Range("D44").Value = "0.002231"; $ws.Range("E44").Value = "6.07%"; $ws.Range("G44").Value = "7"
$ws.Range("D45").Value = "0.009835"; $ws.Range("E45").Value = "3.33%"; $ws.Range("G45").Value = "7"
$ws.Range("D46").Value = "0.00006061"; $ws.Range("E46").Value = "-0.82%"; $ws.Range("G46").Value = "7"
$ws.Range("D47").Value = "0.00000000750"; $ws.Range("E47").Value = "0.17%"; $ws.Range("G47").Value = "7"
$ws.Range("D48").Value = "8.671"; $ws.Range("E48").Value = "217.40%"; $ws.Range("G48").Value = "7"
$ws.Range("D49").Value = "0.002691"; $ws.Range("E49").Value = "34.73%"; $ws.Range("G49").Value = "7"
$ws.Range("D50").Value = "0.00002101"; $ws.Range("E50").Value = "0.17%"; $ws.Range("G50").Value = "7"
$ws.Range("D51").Value = "0.0002001"; $ws.Range("E51").Value = "0.17%"; $ws.Range("G51").Value = "7"
